# Agrego decisiones de diseño
# Applies the design-decision edits to the "Hoja1" worksheet:
#  - reword the "Vinculador" row (now row 25) as "Vinculador Ingresos-Egresos"
#    and move it into the "Otros" scope
#  - rework the Item/ItemEgreso/ItemPresupuesto justification row (row 26)
#  - add a brand new "Otros / Base de datos" row (row 29) about only
#    persisting Argentina, pushing everything below down
#  - add a block of new (currently empty) rows 30-48 below the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 25 — was "Vinculador"; now "Vinculador Ingresos-Egresos" under "Otros"
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Otros"
$ws.Range("B25").Value = "Vinculador Ingresos-Egresos"
$ws.Range("C25").Value = "Encapsulamos el proceso del vinculador en una clase que se encarga de recibir la entidad y los criterios con los que quiere vincular"
$ws.Range("D25").Value = "Hacer un Strategy y una clase para cada metodo de vinculador"

# ---------------------------------------------------------------------
# Row 26 — Item / ItemEgreso / ItemPresupuesto justification reworked
# ---------------------------------------------------------------------
$ws.Range("C26").Value = "Agregamos dos clases que heredan de Item, ItemEgreso e ItemPresupuesto "
$ws.Range("D26").Value = ""
$ws.Range("E26").Font.Underline = $false
$ws.Range("E26").HorizontalAlignment = -4131
$ws.Range("E26").VerticalAlignment = -4160
$ws.Range("E26").WrapText = $true
$ws.Range("E26").Value = "Porque nos parecio mejor tratarlos polimorficamente pero como objetos diferentes. Ademas se nos hacia dificil persistir una unica clase item que se pueda relacionar con dos entidades diferentes "

# ---------------------------------------------------------------------
# Row 29 — new "Base de datos" decision about only persisting Argentina
# (replaces the old "Componente / Brasil" content)
# ---------------------------------------------------------------------
$ws.Range("B29").Value = "Base de datos"
$ws.Range("C29").Value = "En la persistencia de paises solo incluimos a Argentina"
$ws.Range("E29").Value = "Porque sino es un volumen muy grande de datos que no sabemos si son necesarios y relentiza la carga de informacion a la BD"
$ws.Rows.Item(29).RowHeight = 57.6

# ---------------------------------------------------------------------
# Rows 30-48 — new blank rows appended below the table
# ---------------------------------------------------------------------
$ws.Range("A30:C30").HorizontalAlignment = -4131
$ws.Range("A30:C30").VerticalAlignment = -4160
$ws.Range("A30:C30").WrapText = $true

$ws.Range("A31").HorizontalAlignment = -4131
$ws.Range("A31").VerticalAlignment = -4160
$ws.Range("A31").WrapText = $true

$ws.Range("D30:E48").HorizontalAlignment = -4131
$ws.Range("D30:E48").VerticalAlignment = -4160

$ws.Range("B31:E48").HorizontalAlignment = -4131
$ws.Range("B31:E48").VerticalAlignment = -4160

$ws.Range("A32:A48").HorizontalAlignment = -4131
$ws.Range("A32:A48").VerticalAlignment = -4160
